# Use relative performance in benchmark plots
# Switches the "overhead vs reference" (x N) numbers for a "relative
# performance" (x 0.NN, where 1 = reference) view, on both the Linux and
# Windows benchmark sheets, and updates the accompanying charts to match
# (linear scale instead of log10, new titles).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Data sheets (Linux + Windows) - rewrite the "reference multiplier"
#    column (D, text) and the backing numeric column (C) as 1/x formulas.
# ---------------------------------------------------------------------------

# --- Linux sheet ---
$ws1 = $wb.Worksheets.Item("Linux")

$ws1.Range("D3").Value = "× 1 (ref)"
$ws1.Range("C4").Formula = "=1/1.81"
$ws1.Range("D4").Value = "× 0.55"
$ws1.Range("C5").Formula = "=1/75"
$ws1.Range("D5").Value = "× 0.01"

$ws1.Range("D7").Value = "× 1 (ref)"
$ws1.Range("C8").Formula = "=1/1.71"
$ws1.Range("D8").Value = "× 0.58"
$ws1.Range("C9").Formula = "=1/233"
$ws1.Range("D9").Value = "× 0.005"

$ws1.Range("D11").Value = "× 1 (ref)"
$ws1.Range("C12").Formula = "=1/1.18"
$ws1.Range("D12").Value = "× 0.85"
$ws1.Range("C13").Formula = "=1/3.28"
$ws1.Range("D13").Value = "× 0.30"

# --- Windows sheet ---
$ws2 = $wb.Worksheets.Item("Windows")

$ws2.Range("D3").Value = "× 1 (ref)"
$ws2.Range("C4").Formula = "=1/1.84"
$ws2.Range("D4").Value = "× 0.54"
$ws2.Range("C5").Formula = "=1/42"
$ws2.Range("D5").Value = "× 0.02"

$ws2.Range("D7").Value = "× 1 (ref)"
$ws2.Range("C8").Formula = "=1/1.99"
$ws2.Range("D8").Value = "× 0.50"
$ws2.Range("C9").Formula = "=1/161"
$ws2.Range("D9").Value = "× 0.006"

$ws2.Range("D11").Value = "× 1 (ref)"
$ws2.Range("C12").Formula = "=1/1.23"
$ws2.Range("D12").Value = "× 0.81"
$ws2.Range("C13").Formula = "=1/3.7"
$ws2.Range("D13").Value = "× 0.27"

# ---------------------------------------------------------------------------
# 2. Tidy up the now-redundant explicit "applyFont" style on the rand/atoi/
#    raylib reference rows (A3:D5) - it matched the default style already.
# ---------------------------------------------------------------------------

$ws1.Range("A3:C3").ClearFormats()
$ws1.Range("B4:D4").ClearFormats()
$ws1.Range("B5:D5").ClearFormats()

$ws2.Range("A3:C3").ClearFormats()
$ws2.Range("B4:D4").ClearFormats()
$ws2.Range("B5:D5").ClearFormats()

# ---------------------------------------------------------------------------
# 3. Charts - switch both value axes from a log10 scale to a linear one
#    (Linux chart gets an explicit max of 1.2 to leave room for labels),
#    and refresh the chart / axis titles to describe relative performance.
# ---------------------------------------------------------------------------

$co1 = $ws1.ChartObjects().Item(1)
$chart1 = $co1.Chart
$chart1.ChartTitle.Text = "Linux x86_64 (AMD® Ryzen™ 7 4700U)"
$ax1 = $chart1.Axes(2)
$ax1.ScaleType = -4132
$ax1.MaximumScale = 1.2
$ax1.AxisTitle.Text = "Relative performance of each implementation"

$co2 = $ws2.ChartObjects().Item(1)
$chart2 = $co2.Chart
$chart2.ChartTitle.Text = "Windows x86_64 (Intel® Core™ i5-4460)"
$ax2 = $chart2.Axes(2)
$ax2.ScaleType = -4132
$ax2.AxisTitle.Text = "Relative performance of each implementation"

# ---------------------------------------------------------------------------
# 4. View state - Windows sheet becomes the active tab/selection, Linux
#    keeps a parked selection.
# ---------------------------------------------------------------------------

$ws1.Range("P13").Select()
$ws2.Activate()
$ws2.Range("D9").Select()
